$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 5915.875
$ws.Range("I107").Value = 478.57144
$ws.Range("J107").Value = 10144.889
$ws.Range("K107").Value = 478.57144
$ws.Range("L107").Value = 10144.889
$ws.Range("M107").Value = 1441.42856
$ws.Range("N107").Value = -13984.889

$ws.Range("H113").Value = 3003.4375
$ws.Range("I113").Value = 2850.5
$ws.Range("J113").Value = 3258.3333
$ws.Range("K113").Value = 2850.5
$ws.Range("L113").Value = 3258.3333
$ws.Range("M113").Value = 403.5
$ws.Range("N113").Value = -9766.3333

$ws.Range("H116").Value = 3072.7273
$ws.Range("I116").Value = 1600
$ws.Range("J116").Value = 3625
$ws.Range("K116").Value = 1600
$ws.Range("L116").Value = 3625
$ws.Range("M116").Value = 1842
$ws.Range("N116").Value = -10509

$ws.Range("H138").Value = 1745.9877
$ws.Range("I138").Value = 1211.4814
$ws.Range("J138").Value = 2013.2407
$ws.Range("K138").Value = 3634.4442
$ws.Range("L138").Value = 6039.7221
$ws.Range("M138").Value = 1505.5558
$ws.Range("N138").Value = -16319.7221

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 22727884
$ws.Range("I2").Value = 33333694
$ws.Range("J2").Value = 1146.0714
$ws.Range("K2").Value = 33333694
$ws.Range("L2").Value = 1146.0714
$ws.Range("M2").Value = -33333581
$ws.Range("N2").Value = -1372.0714

$ws.Range("H32").Value = 19281.598
$ws.Range("I32").Value = 20796.459
$ws.Range("J32").Value = 13506.1875
$ws.Range("K32").Value = 20796.459
$ws.Range("L32").Value = 13506.1875
$ws.Range("M32").Value = -20509.459
$ws.Range("N32").Value = -14080.1875

$ws.Range("H116").Value = 22727884
$ws.Range("I116").Value = 33333694
$ws.Range("J116").Value = 1146.0714
$ws.Range("K116").Value = 33333694
$ws.Range("L116").Value = 1146.0714
$ws.Range("M116").Value = -33331400
$ws.Range("N116").Value = -5734.0714

$ws.Range("H122").Value = 2072.3635
$ws.Range("I122").Value = 1849.5
$ws.Range("K122").Value = 5548.5
$ws.Range("M122").Value = -3098.5

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 22727884
$ws.Range("I3").Value = 33333694
$ws.Range("J3").Value = 1146.0714
$ws.Range("K3").Value = 33333694
$ws.Range("L3").Value = 1146.0714
$ws.Range("M3").Value = -33333580
$ws.Range("N3").Value = -1374.0714

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3790785.5
$ws.Range("I31").Value = 2149.611
$ws.Range("J31").Value = 20839646
$ws.Range("K31").Value = 2149.611
$ws.Range("L31").Value = 20839646
$ws.Range("M31").Value = -1854.611
$ws.Range("N31").Value = -20840236

$ws.Range("H34").Value = 3790785.5
$ws.Range("I34").Value = 2149.611
$ws.Range("J34").Value = 20839646
$ws.Range("K34").Value = 2149.611
$ws.Range("L34").Value = 20839646
$ws.Range("M34").Value = -1947.611
$ws.Range("N34").Value = -20840050

$ws.Range("H99").Value = 2241.75
$ws.Range("I99").Value = 2007.8334
$ws.Range("J99").Value = 2475.6667
$ws.Range("K99").Value = 2007.8334
$ws.Range("L99").Value = 2475.6667
$ws.Range("M99").Value = -509.8334
$ws.Range("N99").Value = -5471.6667

$ws.Range("H126").Value = 2241.75
$ws.Range("I126").Value = 2007.8334
$ws.Range("J126").Value = 2475.6667
$ws.Range("K126").Value = 6023.5002
$ws.Range("L126").Value = 7427.000100000001
$ws.Range("M126").Value = -3553.5002
$ws.Range("N126").Value = -12367.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 763.97
$ws.Range("J131").Value = 777.88544
$ws.Range("L131").Value = 2333.65632
$ws.Range("N131").Value = -12413.65632

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3073.182
$ws.Range("I80").Value = 3100.8333
$ws.Range("J80").Value = 3040
$ws.Range("K80").Value = 3100.8333
$ws.Range("L80").Value = 3040
$ws.Range("M80").Value = -2102.8333
$ws.Range("N80").Value = -5036

$ws.Range("H83").Value = 3073.182
$ws.Range("I83").Value = 3100.8333
$ws.Range("J83").Value = 3040
$ws.Range("K83").Value = 15504.1665
$ws.Range("L83").Value = 15200
$ws.Range("M83").Value = -10512.1665
$ws.Range("N83").Value = -25184

$ws.Range("H113").Value = 16667466
$ws.Range("J113").Value = 948.5
$ws.Range("L113").Value = 948.5
$ws.Range("N113").Value = -5288.5

$ws.Range("H122").Value = 37039550
$ws.Range("I122").Value = 90911740
$ws.Range("J122").Value = 2419.25
$ws.Range("K122").Value = 272735220
$ws.Range("L122").Value = 7257.75
$ws.Range("M122").Value = -272732770
$ws.Range("N122").Value = -12157.75

$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws.Range("H131").Value = 23442.334
$ws.Range("J131").Value = 23442.334
$ws.Range("L131").Value = 23442.334
$ws.Range("N131").Value = -33522.334

$ws.Range("H132").Value = 43309.92
$ws.Range("I132").Value = 74560.57000000001
$ws.Range("J132").Value = 3536.3635
$ws.Range("K132").Value = 223681.71
$ws.Range("L132").Value = 10609.0905
$ws.Range("M132").Value = -221151.71
$ws.Range("N132").Value = -15669.0905

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 27779166
$ws.Range("I61").Value = 1585
$ws.Range("J61").Value = 83334330
$ws.Range("K61").Value = 1585
$ws.Range("L61").Value = 83334330
$ws.Range("M61").Value = -1383
$ws.Range("N61").Value = -83334734

$ws.Range("H113").Value = 27779166
$ws.Range("I113").Value = 1585
$ws.Range("J113").Value = 83334330
$ws.Range("K113").Value = 1585
$ws.Range("L113").Value = 83334330
$ws.Range("M113").Value = 585
$ws.Range("N113").Value = -83338670

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 7678.1665
$ws.Range("J74").Value = 11133.333
$ws.Range("L74").Value = 11133.333
$ws.Range("N74").Value = -13005.333

$ws.Range("H77").Value = 7678.1665
$ws.Range("J77").Value = 11133.333
$ws.Range("L77").Value = 33399.999
$ws.Range("N77").Value = -42759.999
